$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds decimal price values that are stored as literal text in the
# source sheet (e.g. "250.07"), not as native numbers. Setting .Value directly
# with a numeric-looking string would make Excel coerce the cell to a real
# number (losing trailing zeros / exact text). To keep them as text, flip the
# cell to a text number-format just long enough to assign the value, then put
# the General format back so the cells look just like the rest of the sheet.
$textCells = "D2,D3,D4,D5,D6,D7,D8,D9,D10,D11,D12,D14,D15,D16,D17,D18,D19,D20,D21,D22,D23,D24,D28,D41,D42,D43,D44,D45,D46,D48,D49,D50,D51"
foreach ($r in $textCells.Split(",")) {
    $ws.Range($r).NumberFormat = "@"
}

$ws.Range("D2").Value = "249.11"
$ws.Range("D3").Value = "22.61"
$ws.Range("D4").Value = "5.374"
$ws.Range("D5").Value = "0.05612"
$ws.Range("D6").Value = "3.438"
$ws.Range("D7").Value = "6.353"
$ws.Range("D8").Value = "0.8163"
$ws.Range("D9").Value = "0.9147"
$ws.Range("D10").Value = "0.1418"
$ws.Range("D11").Value = "0.07472"
$ws.Range("D12").Value = "0.03193"
$ws.Range("D14").Value = "0.09327"
$ws.Range("D15").Value = "3.559"
$ws.Range("D16").Value = "0.001594"
$ws.Range("D17").Value = "0.04720"
$ws.Range("D18").Value = "0.0005757"
$ws.Range("E18").Value = "17OneONEWorstin24h"
$ws.Range("D19").Value = "0.006393"
$ws.Range("D20").Value = "0.005000"
$ws.Range("D21").Value = "0.001032"
$ws.Range("D22").Value = "0.0001501"
$ws.Range("D23").Value = "3.727"
$ws.Range("D24").Value = "2.164"
$ws.Range("D28").Value = "0.0002998"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "0.006947"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "0.1065"
$ws.Range("E42").Value = "41BKEXTokenBKK"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "0.003404"
$ws.Range("E43").Value = "42CEJICEJI"
$ws.Range("D44").Value = "0.007560"
$ws.Range("D45").Value = "0.00005569"
$ws.Range("D46").Value = "0.00000000751"
$ws.Range("D48").Value = "0.6748"
$ws.Range("D49").Value = "0.2229"
$ws.Range("D50").Value = "0.00002103"
$ws.Range("D51").Value = "0.01011"

foreach ($r in $textCells.Split(",")) {
    $ws.Range($r).NumberFormat = "General"
}
